# act tablas web jul25
# Adds 2022/2023/2024 data points to the "Data" sheet and updates the
# "Metadata" sheet with a new "actualizacion" entry (Julio 2025) and a
# small tweak to the "cita" text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Data sheet: insert 3 new rows right after the header with the
#    new years (2024, 2023, 2022) and their values, pushing the
#    existing years (2021..2006) down.
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Rows("2:4").Insert()

# Column A stores years as text (shared strings), not numbers. Writing
# a numeric-looking literal via .Value would be auto-detected as a
# number, so instead compute the text with a TEXT() formula and then
# convert the formula results to plain values (Copy + PasteSpecial
# values-only) so no formula nor extra number format is left behind.
$wsData.Cells.Item(2, 1).Formula = "=TEXT(2024,""0"")"
$wsData.Cells.Item(3, 1).Formula = "=TEXT(2023,""0"")"
$wsData.Cells.Item(4, 1).Formula = "=TEXT(2022,""0"")"

$newYearsRange = $wsData.Range("A2:A4")
$newYearsRange.Copy()
$newYearsRange.PasteSpecial(-4163)  # xlPasteValues

$wsData.Cells.Item(2, 2).Value = 22.7
$wsData.Cells.Item(3, 2).Value = 23.4
$wsData.Cells.Item(4, 2).Value = 22.5

# ---------------------------------------------------------------
# 2) Metadata sheet:
#    - the blank key in row 1 becomes a single space (matching the
#      value column), which lets the previously-orphaned empty
#      shared string be dropped.
#    - a new "actualizacion" / "Julio 2025" row is inserted right
#      before the "cita" row.
#    - the "cita" value text gains a trailing newline.
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Cells.Item(1, 1).Value = " "

# Find the "cita" row dynamically and insert the new row above it.
$citaRow = 0
for ($r = 1; $r -le $wsMeta.UsedRange.Rows.Count; $r++) {
    if ($wsMeta.Cells.Item($r, 1).Value2 -eq "cita") {
        $citaRow = $r
        break
    }
}

$wsMeta.Rows("$citaRow`:$citaRow").Insert()

$wsMeta.Cells.Item($citaRow, 1).Value = "actualizacion"
$wsMeta.Cells.Item($citaRow, 2).Value = "Julio 2025"

$citaValueRow = $citaRow + 1
$wsMeta.Cells.Item($citaValueRow, 2).Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE`n"

# Writing text with an embedded newline makes the engine mark the row
# with an explicit custom height; AutoFit() recalculates it back to
# the natural (non-custom) height so no stray row-height attribute is
# left behind.
$wsMeta.Rows($citaValueRow).AutoFit()
